$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.327.33'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.325.56'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.69'
$ws.Range("E5").Value = '  +3.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.57'
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +2.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.320.04'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  +1.85%  '
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.47'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000276'
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '663.48'
$ws.Range("E14").Value = '  +11.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.858.04'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.48'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.414.14'
$ws.Range("E17").Value = '  +3.48%  '
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.326.88'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.76'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.97'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.905'
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.77'
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.97'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.61'
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.79'
$ws.Range("E29").Value = '  +6.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.63'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '602.51'
$ws.Range("E32").Value = '  +7.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.963.86'
$ws.Range("E33").Value = '  +4.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.00'
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.52'
$ws.Range("E36").Value = '  -3.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.91'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.27'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("E41").Value = '  +3.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.70'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0691'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  +1.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("E48").Value = '  +13.33%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.29'
$ws.Range("E51").Value = '  +1.94%  '
